# lab4.0 edit: add note about altera_hostfs only working in debugger mode,
# drop a stray empty "ListParagraph", and shift the _GoBack bookmark to sit
# mid-word as Word's autosave left it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the stray empty ListParagraph (pStyle=ListParagraph,
#    ind left=1440, no numPr, no text) that sits just before the
#    "Once each team's hardware design ..." bullet.
# ---------------------------------------------------------------------
$emptyParas = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r" -and `
        $p.Style.NameLocal -eq "List Paragraph" -and `
        $p.Range.ParagraphFormat.LeftIndent -eq 72) {
        $emptyParas += $p
    }
}
foreach ($p in $emptyParas) {
    $p.Range.Delete()
}

# ---------------------------------------------------------------------
# 2. Append a new explanatory sentence to the paragraph that ends in
#    "... II development tools just as you were doing in the previous
#    labs.", right after that sentence, inside the same paragraph.
#    "altera_hostfs" and "debugger" get the character style "Emphasis",
#    matching how the term is styled elsewhere in the document.
# ---------------------------------------------------------------------
function Insert-Plain($range, [string]$text) {
    $range.Collapse(0) # wdCollapseEnd
    $range.InsertAfter($text)
    return $range
}

function Insert-Emphasis($range, [string]$text) {
    $range.Collapse(0) # wdCollapseEnd
    $range.InsertAfter($text)
    $range.Style = "Emphasis"
    return $range
}

$r = $d.Content
$r.Find.Execute("previous labs.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)

$r = Insert-Plain    $r " Note that the "
$r = Insert-Emphasis $r "altera_hostfs"
$r = Insert-Plain    $r " package only works "
$r = Insert-Plain    $r "when the "
$r = Insert-Emphasis $r "debugger"
$r = Insert-Plain    $r " is being used, not in standard run mode."

# ---------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark so that it splits the word
#    "successfully" into "success" | bookmark | "fully" instead of
#    sitting at the very end of the paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$r2 = $d.Content
$r2.Find.Execute("successfully configure the board", $true, $false, $false, `
                  $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $r2.Start + 7   # length of "success"
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
